# Update the ontodog example worksheet to the latest figure6 / OBCS
# release version: append one new ontology-term row (IRI, label, and
# "include in view" flag) after the existing last data row (row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "http://purl.obolibrary.org/obo/BFO_0000040"
$ws.Range("B22").Value = "material entity"
$ws.Range("C22").Value = "y"

# Reflect the new active selection on the sheet after the edit.
$null = $ws.Range("A22").Select()
